# Apply the updates described by the commit diff:
#  1. Rename the worksheet from "IClientBalance-20240521-085926-" to "Saldo_guide".
#  2. Bump the reference date in column G (rows 2-257) from 45433 to 45434 (2024-05-21 -> 2024-05-22).
#  3. Update row 113's "Saldo Previsto" (D) and "Vl. Total" (H) values from 186 to 223.55.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet/tab.
$ws.Name = "Saldo_guide"

# 2) Shift every reference date in column G (data rows 2 through 257) forward by one day.
$ws.Range("G2:G257").Value = 45434

# 3) Correct the balance figures on row 113.
$ws.Range("D113").Value = 223.55
$ws.Range("H113").Value = 223.55
